$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: updated Price (D) and Volume(1h) (E) columns,
# plus a Dai / WrappedliquidstakedEther2.0 row-order swap (rows 21-22).

$ws.Range('D2').Value = '30.061.87'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '1.903.89'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7468'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.01'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9988'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3078'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.03'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -5.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06911'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08013'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7610'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.78%  '
$ws.Range('D13').Value = '1.917.00'
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.261'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.44'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.196'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.11%  '
$ws.Range('D17').Value = '30.068.17'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.08'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007770'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.36'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.05%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.159.59'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9989'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9982'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.124'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.357'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.15'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.87'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1269'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.051'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -6.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.349'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.530'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.307'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.044'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05353'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.292'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7427'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01944'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.760'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.255'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4468'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.14'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.968'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9987'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8311'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.715'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.44'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.799'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('D49').Value = '2.060.46'
$ws.Range('E49').Value = '  -1.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.72'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1171'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.92%  '
